$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = -12.661
$ws.Range("B7").Value = 5.865
$ws.Range("D7").Value = -7.696000000000001
$ws.Range("A9").Value = -21.881
$ws.Range("D10").Value = -8.337
$ws.Range("B12").Value = 5.782999999999999
$ws.Range("D13").Value = -7.726999999999999
$ws.Range("B14").Value = 6.229
$ws.Range("C15").Value = -12.909
$ws.Range("D16").Value = -8.579000000000001
$ws.Range("A18").Value = -21.993
$ws.Range("A20").Value = -20.931
$ws.Range("D20").Value = -7.874
$ws.Range("D24").Value = -7.547
$ws.Range("B26").Value = 5.962
$ws.Range("A27").Value = -21.334
$ws.Range("B27").Value = 6.006
$ws.Range("B29").Value = 6.169
$ws.Range("C33").Value = -11.314
$ws.Range("A35").Value = -20.067
$ws.Range("C35").Value = -13.013
$ws.Range("B37").Value = 8.270000000000001
$ws.Range("B38").Value = 5.439
$ws.Range("C38").Value = -12.61
$ws.Range("D39").Value = -7.634
$ws.Range("C43").Value = -12.494
$ws.Range("C44").Value = -12.106
$ws.Range("C47").Value = -11.418
$ws.Range("D47").Value = -7.372
$ws.Range("D48").Value = -7.229000000000001
$ws.Range("B51").Value = 5.911
$ws.Range("C51").Value = -11.915
$ws.Range("B52").Value = 5.478
$ws.Range("D52").Value = -7.668000000000001
$ws.Range("B55").Value = 6.394
$ws.Range("D56").Value = -7.825
$ws.Range("C57").Value = -13.667
$ws.Range("C63").Value = -12.277
$ws.Range("A69").Value = -21.701
$ws.Range("B69").Value = 6.616
$ws.Range("B70").Value = 6.201000000000001
$ws.Range("C70").Value = -11.114
$ws.Range("A76").Value = -20.738
$ws.Range("A78").Value = -20.609
$ws.Range("B81").Value = 6.219
$ws.Range("A82").Value = -21.938
$ws.Range("A83").Value = -20.724
$ws.Range("B83").Value = 7.435
$ws.Range("D84").Value = -8.161
$ws.Range("C88").Value = -12.619
$ws.Range("A93").Value = -22.021
$ws.Range("C99").Value = -12.526
$ws.Range("D100").Value = -8.310999999999998
$ws.Range("D101").Value = -7.831
$ws.Range("B102").Value = 7.402999999999999
